{"js": "// Remove the trailing \"Ver no Jupiter ...\" line, the \"\u00a9 2020 ...\" credit\n// line, and the blank paragraph that separates them from the preceding\n// \"Requisitos\" content \u2014 mirrors the upstream Jekyll rebuild that dropped\n// the site-chrome boilerplate from the scraped page.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nconst jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst copyrightPrefix = \"\u00a9 2020\";\n\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\n\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (jupiterIndex === -1 && t.trim() === jupiterText) {\n    jupiterIndex = i;\n  } else if (t.trim().indexOf(copyrightPrefix) === 0) {\n    copyrightIndex = i;\n  }\n}\n\nif (jupiterIndex === -1 || copyrightIndex === -1) {\n  throw new Error(\"Could not locate the Jupiter/copyright paragraphs to remove.\");\n}\n\n// The blank separator paragraph (if present) directly precedes the\n// \"Ver no Jupiter ...\" paragraph.\nlet blankIndex = -1;\nif (jupiterIndex > 0 && items[jupiterIndex - 1].text.trim() === \"\") {\n  blankIndex = jupiterIndex - 1;\n}\n\n// Delete from the bottom up so earlier indices stay valid.\nconst toDelete = [copyrightIndex, jupiterIndex];\nif (blankIndex !== -1) {\n  toDelete.push(blankIndex);\n}\ntoDelete.sort((a, b) => b - a);\n\nfor (const idx of toDelete) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" line, the \"(c) 2020 ...\" credit\n# line, and the blank paragraph that separates them from the preceding\n# \"Requisitos\" content -- mirrors the upstream Jekyll rebuild that dropped\n# the site-chrome boilerplate from the scraped page.\n$d = $word.ActiveDocument\n\n$jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightPrefix = [string][char]0x00A9 + \" 2020\"\n\n$jupiterIndex = -1\n$copyrightIndex = -1\n$count = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r`a\", \"`r\", \"`a\")\n    $t = $t.Trim()\n    if ($jupiterIndex -eq -1 -and $t -eq $jupiterText) {\n        $jupiterIndex = $i\n    } elseif ($t.StartsWith($copyrightPrefix)) {\n        $copyrightIndex = $i\n    }\n}\n\nif ($jupiterIndex -eq -1 -or $copyrightIndex -eq -1) {\n    throw \"Could not locate the Jupiter/copyright paragraphs to remove.\"\n}\n\n$blankIndex = -1\nif ($jupiterIndex -gt 1) {\n    $prevText = $d.Paragraphs.Item($jupiterIndex - 1).Range.Text.TrimEnd(\"`r`a\", \"`r\", \"`a\").Trim()\n    if ($prevText -eq \"\") {\n        $blankIndex = $jupiterIndex - 1\n    }\n}\n\n$toDelete = @($copyrightIndex, $jupiterIndex)\nif ($blankIndex -ne -1) {\n    $toDelete += $blankIndex\n}\n$toDelete = $toDelete | Sort-Object -Descending\n\nforeach ($idx in $toDelete) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
